# CRYPTO61675.xlsx — "updated headers of all HMRC examples"
#
# The meaningful edit: three table/header columns get " in GBP" appended to
# make it clear the money columns are denominated in GBP:
#   "Buy Value"  -> "Buy Value in GBP"
#   "Sell Value" -> "Sell Value in GBP"
#   "Fee Value"  -> "Fee Value in GBP"
#
# Renaming the header cell text (rather than poking ListColumns(i).Name,
# which is not wired up for writes here) updates both the shared-strings
# table and the table1.xml column definitions, since the table's header
# row is backed directly by these worksheet cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Buy Value in GBP"
$ws.Range("G1").Value = "Sell Value in GBP"
$ws.Range("J1").Value = "Fee Value in GBP"

# Column widths were refreshed (auto best-fit) to accommodate the longer
# header text. Set the resulting widths directly (COM ColumnWidth is in
# character units; the engine adds ~5/6 of a character for cell padding
# when it serialises the sheet's <col> widths).
$widths = @{
    1  = 6.5
    2  = 11.16666666666667
    3  = 8.5
    4  = 14.66666666666667
    5  = 11
    6  = 8.333333333333334
    7  = 14.5
    8  = 11
    9  = 8.333333333333334
    10 = 14.5
    11 = 11.833333333333334
    12 = 17.333333333333332
    13 = 255
}
foreach ($colIndex in $widths.Keys) {
    $ws.Columns.Item($colIndex).ColumnWidth = $widths[$colIndex]
}

# Cursor/selection moved as part of the editing session.
$ws.Range("H19").Select()
